$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column BB (col 54) -------------------------------------------------
# Rows 1-71: identical to column BA (col 53), so copy format + values across.
$ws.Range("BA1:BA71").Copy()
$ws.Range("BB1:BB71").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("BA1:BA71").Copy()
$ws.Range("BB1:BB71").PasteSpecial(-4163)   # xlPasteValues
$excel.CutCopyMode = 0

# Header date for the new quarter column differs from BA1.
$ws.Cells.Item(1, 54).Value = 45986

# Rows 72-82: new forecast values (format copied from BA column, same style).
$ws.Range("BA72:BA82").Copy()
$ws.Range("BB72:BB82").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Cells.Item(72, 54).Value = -0.2099029780610664
$ws.Cells.Item(73, 54).Value = 0
$ws.Cells.Item(74, 54).Value = -0.0492253650248415
$ws.Cells.Item(75, 54).Value = -0.0492253650248415
$ws.Cells.Item(76, 54).Value = -0.0492253650248415
$ws.Cells.Item(77, 54).Value = -0.0492253650248415
$ws.Cells.Item(78, 54).Value = -0.0492253650248415
$ws.Cells.Item(79, 54).Value = -0.0492253650248415
$ws.Cells.Item(80, 54).Value = -0.0492253650248415
$ws.Cells.Item(81, 54).Value = -0.0492253650248415
$ws.Cells.Item(82, 54).Value = -0.0492253650248415

# --- New row 83 --------------------------------------------------------------
# A83 carries the next quarter's date, styled like the rest of column A.
$ws.Range("A82").Copy()
$ws.Range("A83").PasteSpecial(-4122)        # xlPasteFormats
$excel.CutCopyMode = 0
$ws.Cells.Item(83, 1).Value = 46934

# BB83 new forecast value (plain, unstyled, like the rest of column BB).
$ws.Cells.Item(83, 54).Value = -0.0492253650248415
